# Populate Price (H), Stock (I) and Status (J) columns for the supplier
# rows (5-44) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(5, 3000, 25, "Available"),
    @(6, 5000, 50, "Available"),
    @(7, 4000, 0, "Out of Stock"),
    @(8, 3500, 50, "Available"),
    @(9, 4500, 0, "Out of Stock"),
    @(10, 3000, 75, "Available"),
    @(11, 5000, 25, "Unavailable"),
    @(12, 3500, 50, "Available"),
    @(13, 4500, 100, "Available"),
    @(14, 4000, 0, "Out of Stock"),
    @(15, 3500, 25, "Available"),
    @(16, 5000, 75, "Available"),
    @(17, 4000, 0, "Out of Stock"),
    @(18, 3500, 100, "Available"),
    @(19, 3500, 50, "Available"),
    @(20, 5000, 25, "Unavailable"),
    @(21, 4500, 100, "Available"),
    @(22, 3000, 75, "Available"),
    @(23, 4500, 50, "Available"),
    @(24, 4000, 25, "Unavailable"),
    @(25, 3000, 50, "Available"),
    @(26, 3500, 0, "Out of Stock"),
    @(27, 5000, 100, "Available"),
    @(28, 3500, 25, "Available"),
    @(29, 3000, 100, "Available"),
    @(30, 4500, 75, "Available"),
    @(31, 5000, 25, "Unavailable"),
    @(32, 3000, 50, "Available"),
    @(33, 4500, 100, "Available"),
    @(34, 4000, 75, "Available"),
    @(35, 3500, 0, "Out of Stock"),
    @(36, 4500, 25, "Unavailable"),
    @(37, 5000, 100, "Available"),
    @(38, 3500, 75, "Available"),
    @(39, 4000, 50, "Available"),
    @(40, 3000, 25, "Unavailable"),
    @(41, 4500, 100, "Available"),
    @(42, 4000, 0, "Out of Stock"),
    @(43, 3000, 75, "Available"),
    @(44, 5000, 25, "Available")
)

# The shared-strings table records each unique string in first-seen
# order. To reproduce the exact table order of the source workbook
# ("Out of Stock", "Unavailable", "Available"), prime the Status column
# with one cell of each value, in that order, before filling the rest.
$ws.Cells.Item(7, 10).Value = "Out of Stock"
$ws.Cells.Item(11, 10).Value = "Unavailable"
$ws.Cells.Item(5, 10).Value = "Available"

foreach ($row in $data) {
    $r = $row[0]
    $price = $row[1]
    $stock = $row[2]
    $status = $row[3]

    $ws.Cells.Item($r, 8).Value = $price
    $ws.Cells.Item($r, 9).Value = $stock
    $ws.Cells.Item($r, 10).Value = $status
}

# Reset the selection (also clears the saved scrolled-to top-left cell).
$ws.Range("J5").Select()
